# The data rows (2..23) got reshuffled: every row now carries the
# "record" (Fecha, Calidad, Volumen, Precio min/max/promedio, Unidad de
# comercialización, Precio $/Kg, Kg o Unidades) that used to live in a
# different row of the same sheet, while the other columns
# (Mercado ID/Mercado/Región/Codreg/Categoría ID/Categoría/Variedad/
# Origen/Clasificación) stay identical because they were constant for
# every row anyway.
#
# new row -> old row the record comes from
$map = @{
    2  = 16
    3  = 13
    4  = 7
    5  = 20
    6  = 21
    7  = 17
    8  = 9
    9  = 19
    10 = 2
    11 = 6
    12 = 23
    13 = 22
    14 = 8
    15 = 3
    16 = 10
    17 = 15
    18 = 12
    19 = 4
    20 = 18
    21 = 5
    22 = 14
    23 = 11
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one "record" when rows are reshuffled.
$cols = @(4, 9, 10, 11, 12, 13, 14, 16, 17)   # D, I, J, K, L, M, N, P, Q

# Snapshot every original value first, because the mapping is a full
# permutation and rows are both sources and destinations.
$snapshot = @{}
for ($r = 2; $r -le 23; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Now write back according to the mapping.
for ($r = 2; $r -le 23; $r++) {
    $src = $map[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}
